$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reference style with no explicit formatting, used to restore default (no "s" attr) style
$defaultStyle = $ws.Range("B2").Style

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '69.213.16'
$ws.Range("D2").Style = $defaultStyle
$ws.Range("E2").Value = '  +0.42%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.757.01'
$ws.Range("D3").Style = $defaultStyle
$ws.Range("E3").Value = '  +0.62%  '
$ws.Range("E4").Value = '  -0.04%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '602.57'
$ws.Range("D5").Style = $defaultStyle
$ws.Range("E5").Value = '  +0.13%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '167.36'
$ws.Range("D6").Style = $defaultStyle
$ws.Range("E6").Value = '  -0.86%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '3.754.92'
$ws.Range("D7").Style = $defaultStyle
$ws.Range("E7").Value = '  +0.64%  '
$ws.Range("E8").Value = '  -0.05%  '
$ws.Range("E9").Value = '  +1.21%  '
$ws.Range("E10").Value = '  +3.46%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '6.40'
$ws.Range("D11").Style = $defaultStyle
$ws.Range("E11").Value = '  +2.12%  '
$ws.Range("E12").Value = '  -0.30%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '38.09'
$ws.Range("D13").Style = $defaultStyle
$ws.Range("E13").Value = '  -0.22%  '
$ws.Range("E14").Value = '  +2.06%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '4.382.18'
$ws.Range("D15").Style = $defaultStyle
$ws.Range("E15").Value = '  +0.49%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.744.34'
$ws.Range("D16").Style = $defaultStyle
$ws.Range("E16").Value = '  +0.23%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '69.210.08'
$ws.Range("D17").Style = $defaultStyle
$ws.Range("E17").Value = '  +0.41%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '7.37'
$ws.Range("D18").Style = $defaultStyle
$ws.Range("E18").Value = '  +1.67%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '17.38'
$ws.Range("D19").Style = $defaultStyle
$ws.Range("E19").Value = '  +0.94%  '
$ws.Range("E20").Value = '  -1.46%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '11.44'
$ws.Range("D21").Style = $defaultStyle
$ws.Range("E21").Value = '  +20.27%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '494.26'
$ws.Range("D22").Style = $defaultStyle
$ws.Range("E22").Value = '  -0.48%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.730'
$ws.Range("D23").Style = $defaultStyle
$ws.Range("E23").Value = '  +0.95%  '
$ws.Range("E24").Value = '  +7.96%  '
$ws.Range("E25").Value = '  +0.08%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.31'
$ws.Range("D26").Style = $defaultStyle
$ws.Range("E26").Value = '  -0.22%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '12.35'
$ws.Range("D27").Style = $defaultStyle
$ws.Range("E27").Value = '  +0.67%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '10.14'
$ws.Range("D28").Style = $defaultStyle
$ws.Range("E28").Value = '  +0.23%  '
$ws.Range("E29").Value = '  -0.05%  '
$ws.Range("E30").Value = '  +1.87%  '
$ws.Range("E31").Value = '  +2.32%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '8.13'
$ws.Range("D32").Style = $defaultStyle
$ws.Range("E32").Value = '  +2.18%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '31.64'
$ws.Range("D33").Style = $defaultStyle
$ws.Range("E33").Value = '  +0.06%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.898.27'
$ws.Range("D34").Style = $defaultStyle
$ws.Range("E34").Value = '  +0.59%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '3.691.73'
$ws.Range("D35").Style = $defaultStyle
$ws.Range("E35").Value = '  +0.58%  '
$ws.Range("E36").Value = '  -0.59%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.999'
$ws.Range("D37").Style = $defaultStyle
$ws.Range("E37").Value = '  -0.12%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '6.00'
$ws.Range("D38").Style = $defaultStyle
$ws.Range("E38").Value = '  +3.63%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.02'
$ws.Range("D39").Style = $defaultStyle
$ws.Range("E39").Value = '  +0.53%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.136'
$ws.Range("D40").Style = $defaultStyle
$ws.Range("E40").Value = '  +1.75%  '
$ws.Range("E41").Value = '  +0.62%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '3.02'
$ws.Range("D42").Style = $defaultStyle
$ws.Range("E42").Value = '  +5.53%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '430.89'
$ws.Range("D43").Style = $defaultStyle
$ws.Range("E43").Value = '  -1.19%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '48.76'
$ws.Range("D44").Style = $defaultStyle
$ws.Range("E44").Value = '  -0.50%  '
$ws.Range("E45").Value = '  +0.03%  '
$ws.Range("E46").Value = '  +1.22%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '40.40'
$ws.Range("D48").Style = $defaultStyle
$ws.Range("E48").Value = '  -0.07%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '140.94'
$ws.Range("D49").Style = $defaultStyle
$ws.Range("E49").Value = '  -1.73%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.798.61'
$ws.Range("D50").Style = $defaultStyle
$ws.Range("E50").Value = '  +1.68%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0353'
$ws.Range("D51").Style = $defaultStyle
$ws.Range("E51").Value = '  +0.38%  '
